$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2285
$ws.Range("F3").Value = 357
$ws.Range("F4").Value = 185
$ws.Range("F5").Value = 294
$ws.Range("F6").Value = 458
$ws.Range("F8").Value = 716
$ws.Range("F9").Value = 521
$ws.Range("F10").Value = 684
$ws.Range("F11").Value = 373
$ws.Range("F12").Value = 66
$ws.Range("F13").Value = 365
$ws.Range("F14").Value = 977
$ws.Range("F15").Value = 7784
$ws.Range("F16").Value = 219
$ws.Range("F17").Value = 20
$ws.Range("F18").Value = 43
$ws.Range("F19").Value = 254
$ws.Range("F20").Value = 145
$ws.Range("F21").Value = 114
$ws.Range("F23").Value = 107
$ws.Range("F25").Value = 270
$ws.Range("F26").Value = 112

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 211
$ws.Range("F8").Value = 2899
$ws.Range("F10").Value = 26
$ws.Range("F16").Value = 2592

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 52
$ws.Range("F4").Value = 427
$ws.Range("F5").Value = 177

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 52
$ws.Range("F6").Value = 2285
$ws.Range("F7").Value = 427
$ws.Range("F8").Value = 357
$ws.Range("F9").Value = 185
$ws.Range("F10").Value = 294
$ws.Range("F11").Value = 458
$ws.Range("F16").Value = 177
$ws.Range("F17").Value = 716
$ws.Range("F18").Value = 521
$ws.Range("F19").Value = 684
$ws.Range("F20").Value = 373
$ws.Range("F21").Value = 66
$ws.Range("F22").Value = 365
$ws.Range("F23").Value = 977
$ws.Range("F24").Value = 7788
$ws.Range("F25").Value = 211
$ws.Range("F26").Value = 2899
$ws.Range("F28").Value = 26
$ws.Range("F30").Value = 219
$ws.Range("F31").Value = 20
$ws.Range("F32").Value = 43
$ws.Range("F35").Value = 254
$ws.Range("F36").Value = 145
$ws.Range("F37").Value = 114
$ws.Range("F41").Value = 107
$ws.Range("F43").Value = 270
$ws.Range("F44").Value = 112
$ws.Range("F45").Value = 2592
